$wb = $excel.ActiveWorkbook

# --- "Schedule" sheet: update rows 3-5, then remove row 6 entirely ---
$sched = $wb.Worksheets.Item("Schedule")

$sched.Range("A3").Value = 46065.375
$sched.Range("C3").Value = 6.5
$sched.Range("D3").Value = 24.57
$sched.Range("E3").Value = 1160.709849
$sched.Range("F3").Value = 47.24093809523809
$sched.Range("A4").Value = 46065.875
$sched.Range("B4").Value = 46066.14583333334
$sched.Range("C4").Value = 6.5
$sched.Range("D4").Value = 24.57
$sched.Range("E4").Value = 833.799759
$sched.Range("F4").Value = 33.93568412698413
$sched.Range("A5").Value = 46066.3125
$sched.Range("C5").Value = 8.5
$sched.Range("D5").Value = 32.13
$sched.Range("E5").Value = 951.1251750000004
$sched.Range("F5").Value = 29.60240196078433

# Row 6 (the old last row) is dropped from this run; deleting the whole row
# shifts the dimension from A1:F6 down to A1:F5.
$sched.Rows.Item(6).Delete()

# --- "Detailed" sheet: refresh Price / Type / Pump_Status for this run ---
$det = $wb.Worksheets.Item("Detailed")

$det.Range("E20").Value = "ON"
$det.Range("E21").Value = "ON"
$det.Range("B36").Value = 75.00112
$det.Range("B37").Value = 92.04031000000001
$det.Range("B38").Value = 57.31
$det.Range("B39").Value = 57.31
$det.Range("B40").Value = 36.25
$det.Range("C40").Value = "historical"
$det.Range("B41").Value = 115
$det.Range("C41").Value = "historical"
$det.Range("B42").Value = 299.99
$det.Range("C42").Value = "historical"
$det.Range("E42").Value = "OFF"
$det.Range("B43").Value = 139.51244
$det.Range("C43").Value = "historical"
$det.Range("E43").Value = "OFF"
$det.Range("B44").Value = 73.43344999999999
$det.Range("C44").Value = "historical"
$det.Range("B45").Value = 75.71758
$det.Range("C45").Value = "historical"
$det.Range("B46").Value = 84.79000000000001
$det.Range("C46").Value = "historical"
$det.Range("B47").Value = 79.95038
$det.Range("C47").Value = "historical"
$det.Range("B48").Value = 73.20012
$det.Range("C48").Value = "historical"
$det.Range("B49").Value = 57.31
$det.Range("C49").Value = "historical"
$det.Range("B50").Value = 57.68
$det.Range("E50").Value = "ON"
$det.Range("B51").Value = 66.54771
$det.Range("E51").Value = "ON"
$det.Range("B52").Value = 57.31
$det.Range("E52").Value = "ON"
$det.Range("B53").Value = 57.31
$det.Range("E53").Value = "ON"
$det.Range("B54").Value = 57.31
$det.Range("E54").Value = "ON"
$det.Range("B55").Value = 57.31
$det.Range("E55").Value = "ON"
$det.Range("B56").Value = 57.31
$det.Range("E56").Value = "ON"
$det.Range("B57").Value = 57.31
$det.Range("B58").Value = 79.95038
$det.Range("B59").Value = 64.89
$det.Range("B60").Value = 64.89
$det.Range("B61").Value = 71.1951
$det.Range("B62").Value = 83.48781
$det.Range("B63").Value = 84.79000000000001
$det.Range("B64").Value = 73.19
$det.Range("E65").Value = "ON"
$det.Range("B66").Value = 57.0595
$det.Range("B69").Value = 57.06007
$det.Range("B70").Value = 60.57149
$det.Range("B71").Value = 58.30307
$det.Range("B72").Value = 57.8781
$det.Range("B75").Value = 57.06007
$det.Range("B78").Value = 56.98
$det.Range("B79").Value = 57.06007
$det.Range("B80").Value = 57.06007
$det.Range("B83").Value = 56.98
$det.Range("B84").Value = 53.90789
$det.Range("B85").Value = 53.83945
$det.Range("B86").Value = 48.11085
$det.Range("B87").Value = 39.58292
$det.Range("B88").Value = 64.89
$det.Range("B89").Value = 68.6712
$det.Range("B90").Value = 66.42968
$det.Range("E90").Value = "OFF"
$det.Range("B91").Value = 68.35113
$det.Range("E91").Value = "OFF"
$det.Range("B92").Value = 64.89
$det.Range("E92").Value = "OFF"
$det.Range("B93").Value = 64.89
$det.Range("E93").Value = "OFF"
$det.Range("B94").Value = 64.89
$det.Range("E94").Value = "OFF"
$det.Range("B95").Value = 73.19
$det.Range("E95").Value = "OFF"
$det.Range("B96").Value = 67.81603
$det.Range("E96").Value = "OFF"
$det.Range("B97").Value = 69.28548000000001
$det.Range("E97").Value = "OFF"
